# Update the answer values in the division worksheet table.
# The document contains a single 20-row x 5-column table where the
# five "data" rows (1, 5, 9, 13, 17) hold the division problems/answers
# and the other rows are blank spacer rows.
#
# We address each target cell directly by (row, column) and overwrite its
# Range.Text. This preserves run formatting (font/size) already present in
# the cell and avoids any ambiguity that a global Find/Replace could run
# into when a "new" value happens to equal another cell's "old" value
# (e.g. "396÷4=99, 0" is both a target value for one cell and the
# pre-existing text of another cell).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "960÷2=480, 0"
$t.Cell(1, 2).Range.Text  = "743÷9=82, 5"
$t.Cell(1, 3).Range.Text  = "784÷7=112, 0"
$t.Cell(1, 4).Range.Text  = "890÷8=111, 2"
$t.Cell(1, 5).Range.Text  = "312÷7=44, 4"

$t.Cell(5, 1).Range.Text  = "657÷4=164, 1"
$t.Cell(5, 2).Range.Text  = "631÷4=157, 3"
$t.Cell(5, 3).Range.Text  = "880÷5=176, 0"
$t.Cell(5, 4).Range.Text  = "763÷8=95, 3"
$t.Cell(5, 5).Range.Text  = "755÷4=188, 3"

$t.Cell(9, 1).Range.Text  = "219÷8=27, 3"
$t.Cell(9, 2).Range.Text  = "396÷4=99, 0"
$t.Cell(9, 3).Range.Text  = "245÷5=49, 0"
$t.Cell(9, 4).Range.Text  = "569÷8=71, 1"
$t.Cell(9, 5).Range.Text  = "882÷5=176, 2"

$t.Cell(13, 1).Range.Text = "574÷3=191, 1"
$t.Cell(13, 2).Range.Text = "550÷4=137, 2"
$t.Cell(13, 3).Range.Text = "223÷4=55, 3"
$t.Cell(13, 4).Range.Text = "553÷7=79, 0"
$t.Cell(13, 5).Range.Text = "980÷9=108, 8"

$t.Cell(17, 1).Range.Text = "331÷8=41, 3"
$t.Cell(17, 2).Range.Text = "293÷3=97, 2"
$t.Cell(17, 3).Range.Text = "685÷9=76, 1"
$t.Cell(17, 4).Range.Text = "139÷5=27, 4"
$t.Cell(17, 5).Range.Text = "590÷8=73, 6"
